$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update number format of B8:B13 (from General/centered to 0.00/centered -> style s="7")
$ws.Range("B8:B13").NumberFormat = "0.00"

# New row 14: date, hours, description
$ws.Cells.Item(14, 1).Value = 42893
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Cells.Item(14, 2).Value = 1.5
$ws.Cells.Item(14, 2).NumberFormat = "0.00"
$ws.Cells.Item(14, 3).Value = "Nové obrázky na index+parallax, doplnění textů"

# Rows 15-36: just an empty, formatted B cell each (placeholders for future entries)
$ws.Range("B15:B36").NumberFormat = "0.00"

# Restore the active cell selection to C14 (was C16 before the edit)
[void]$ws.Range("C14").Select()
